# Adds an empty run (<w:r><w:t/></w:r>) to a set of paragraphs that were
# previously fully empty (<w:p/> or a <w:pPr>-only paragraph). This mirrors
# the upstream commit "Added specification for requirement API upate." which
# introduced placeholder empty runs in several blank paragraphs across the
# main body, the header and the footer.

$d = $word.ActiveDocument

# Paragraph indices (1-based, Word's Paragraphs collection) of the blank
# paragraphs in the document body that need an empty run appended.
#   1,2,3  -> blank centered title-page paragraphs (before "Python4Capella
#             Simplified Metamodel")
#   5      -> blank centered paragraph after the title
#   30     -> blank paragraph right after the first table
#   95,115 -> blank underlined paragraphs ("u val=single") before the
#             "[CDB] Capella Light Metamodel:" / similar headings
$targets = @(1, 2, 3, 5, 30, 95, 115)

foreach ($i in $targets) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.InsertAfter("")
}

# Header: append an empty run to the final (otherwise empty) paragraph.
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrRange = $hdr.Range
$hdrRange.Collapse(0)
$hdrRange.InsertAfter("")

# Footer: same treatment for its trailing empty paragraph.
$ftr = $sec.Footers.Item(1)
$ftrRange = $ftr.Range
$ftrRange.Collapse(0)
$ftrRange.InsertAfter("")
